# Add "limndiat" as a newly identified missing variable, just before the
# existing "limndiaz" row (row 184), shifting limndiaz and everything below
# it down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 184 (pushes limndiaz, etc. down to row 185+).
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the limndiat data.
$ws.Range("B184").Value = "limndiat"
$ws.Range("H184").Value = "Identified in the shaconemo (r274) ocnBgchem ping file: LDnutSFC"
$ws.Range("I184").Value = "Thomas Reerink"
$ws.Range("L184").Value = "ocnBgchem"
$ws.Range("M184").Value = 1
$ws.Range("N184").Value = "P1 (1) nitrogen_growth_limitation_of_diatoms : Diatoms are phytoplankton with an external skeleton made of silica. Phytoplankton are algae that grow where there is sufficient light to support photosynthesis. ""Nitrogen growth limitation"" means the ratio of the growth rate of a species population in the environment (where there is a finite availability of nitrogen) to the theoretical growth rate if there were no such limit on nitrogen availability."
